# Commit: "test finished, try semi supervise model"
# Adds a new set of Precision/Recall/F1/support test-run columns (G:K) for
# the TextRNN_Attention vs Bert comparison block (rows 16-20) on the
# "Comment" sheet, and moves the active selection to H24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comment")

# --- New block: rows 16-20, columns G:K -------------------------------
# Row 16 (label 0)
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0.8654
$ws.Range("I16").Value = 0.9474
$ws.Range("J16").Value = 0.9045
$ws.Range("K16").Value = 95

# Row 17 (label 1)
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 0.9479
$ws.Range("I17").Value = 0.8667
$ws.Range("J17").Value = 0.9055
$ws.Range("K17").Value = 105

# Row 18 (accuracy) - only overall F1/support populated
$ws.Range("G18").Value = "accuracy"
$ws.Range("J18").Value = 0.905
$ws.Range("K18").Value = 200

# Row 19 (macro avg)
$ws.Range("G19").Value = "macro avg"
$ws.Range("H19").Value = 0.9067
$ws.Range("I19").Value = 0.907
$ws.Range("J19").Value = 0.905
$ws.Range("K19").Value = 200

# Row 20 (weighted avg)
$ws.Range("G20").Value = "weighted avg"
$ws.Range("H20").Value = 0.9087
$ws.Range("I20").Value = 0.905
$ws.Range("J20").Value = 0.905
$ws.Range("K20").Value = 200

# --- Move the active selection to H24 ---------------------------------
$ws.Activate()
$ws.Range("H24").Select()
